# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect the latest scrape, per the commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Changes on the "展览" sheet (row -> new value for column F)
$exhibitionChanges = @{
    7  = 1290
    12 = 204
    22 = 141
    23 = 1149
    24 = 389
    26 = 909
    27 = 1188
    33 = 91
    34 = 578
    36 = 1611
    38 = 1650
    41 = 819
    42 = 17
    43 = 769
    44 = 754
    45 = 956
    46 = 413
    47 = 3293
}

foreach ($row in $exhibitionChanges.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionChanges[$row]
}

# Changes on the "全部类型" sheet (row -> new value for column F)
$allTypesChanges = @{
    9  = 1290
    14 = 204
    23 = 389
    27 = 1188
    32 = 91
    35 = 578
    37 = 1611
    40 = 1650
    41 = 819
    42 = 769
    43 = 754
    44 = 956
    45 = 413
    48 = 3293
}

foreach ($row in $allTypesChanges.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesChanges[$row]
}
